$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "Blackbox Ai" website entry in row 17
$ws.Range("A17").Value = "Blackbox Ai"
$ws.Range("B17").Value = "https://www.blackbox.ai/"
$ws.Range("C17").Value = "https://s13.gifyu.com/images/S08Tb.png"
$ws.Range("D17").Value = "Enhance your coding experience with our powerful AI assistant tool and AI Code Chat - the ultimate companion for programmers. With intelligent suggestions, code completion, and error detection, our tool accelerates software development. With Blackbox you get coding support, AI, intelligent automation, error detection, code completion, and programming efficiency to optimize productivity and streamline your coding workflow."
$ws.Range("E17").Value = "For Developers"
